$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0.5290000438690186
$ws.Range("B4").Value = 0.26200008392333984
$ws.Range("B5").Value = 1.312000036239624
$ws.Range("B6").Value = 0.5769999027252197
$ws.Range("B7").Value = 0.7939999103546143
$ws.Range("B8").Value = 30.063000202178955
$ws.Range("B9").Value = 7.355000019073486
$ws.Range("B10").Value = 30.062000036239624
$ws.Range("B11").Value = 30.158999919891357
$ws.Range("B12").Value = 30.068000078201294
$ws.Range("B13").Value = 30.103000164031982
$ws.Range("E13").Value = 250.72464559656783
$ws.Range("F13").Value = 0.8033584626809539
$ws.Range("B14").Value = 30.200000047683716
